$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Add 29 May 2020 data (row 79)
$ws.Cells.Item(79, 1).Value = 43980
$ws.Cells.Item(79, 2).Value = 36155
$ws.Cells.Item(79, 3).Value = 1141
$ws.Cells.Item(79, 4).Value = 28
$ws.Cells.Item(79, 5).Value = 1594

# Add 30 May 2020 data (row 80)
$ws.Cells.Item(80, 1).Value = 43981
$ws.Cells.Item(80, 2).Value = 39230
$ws.Cells.Item(80, 3).Value = 983
$ws.Cells.Item(80, 4).Value = 26
$ws.Cells.Item(80, 5).Value = 1021

# Extend the table range to include the new rows
$table = $ws.ListObjects.Item("Table3")
$table.Resize($ws.Range("A1:E80"))

# Update selection to match the new last cell, mirroring the saved view state
$ws.Range("E80").Select()
